# DeepBench_NV_TitanXp.xlsx — "Clean up output width and height for convolutions"
#
# Sheet "Results - Mixed Precision" (the active sheet) has, for the convolution
# rows 111-217, two helper columns:
#   O = output width / height, originally computed straight off D/I/K/M
#   P = output width / height, originally computed straight off C/H/J/L
# The commit replaces the naive "(X - Y + 1 + 2*Z)/W" formulas with the correct
# "1 + ROUNDDOWN((X - Y + 2*Z)/W, 0)" formula AND swaps which set of columns O
# and P each read from (O now uses C/H/J/L, P now uses D/I/K/M). Q (TERAFLOPS)
# is left as the same formula but naturally recomputes off the corrected O/P.
#
# Row 111 keeps its own (non-shared) formulas in O/P, matching how the rest of
# the column block had already been split into shared-formula groups at
# O112:O175 / P112:P175 and O176:O217 / P176:P217; Q gets its own shared
# groups split at Q111:Q142, Q143:Q174, Q175:Q206, Q207:Q217 (last row that
# actually holds data in this block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 111: standalone formulas (not part of a fill-down block) ---
$ws.Range("O111").Formula = '=1+ROUNDDOWN((($C111-$H111+2*$J111)/$L111),0)'
$ws.Range("P111").Formula = '=1+ROUNDDOWN((($D111-$I111+2*$K111)/$M111),0)'

# --- O / P fill-down blocks (these become shared formulas across the range) ---
$ws.Range("O112:O175").Formula = '=1+ROUNDDOWN((($C112-$H112+2*$J112)/$L112),0)'
$ws.Range("P112:P175").Formula = '=1+ROUNDDOWN((($D112-$I112+2*$K112)/$M112),0)'

$ws.Range("O176:O217").Formula = '=1+ROUNDDOWN((($C176-$H176+2*$J176)/$L176),0)'
$ws.Range("P176:P217").Formula = '=1+ROUNDDOWN((($D176-$I176+2*$K176)/$M176),0)'

# --- Q (TERAFLOPS) recompute — same formula shape, new shared-formula blocks ---
$ws.Range("Q111:Q142").Formula = '=(2*O111*P111*E111*F111*G111*I111*H111)/(N111/1000)/10^12'
$ws.Range("Q143:Q174").Formula = '=(2*O143*P143*E143*F143*G143*I143*H143)/(N143/1000)/10^12'
$ws.Range("Q175:Q206").Formula = '=(2*O175*P175*E175*F175*G175*I175*H175)/(N175/1000)/10^12'
$ws.Range("Q207:Q217").Formula = '=(2*O207*P207*E207*F207*G207*I207*H207)/(N207/1000)/10^12'

# --- Selection moved down one row (A5 -> A6) ---
$ws.Range("A6").Select()
